# Improving Daftar Jual Saya and Search Product feature
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update product names (shared strings change value but keep their relative
# row mapping swapped: what was in A2/A4 is now swapped, and the text itself
# is replaced with the new product names)
$ws.Range("A2").Value = "zenbook"
$ws.Range("A3").Value = "flanel"
$ws.Range("A4").Value = "hrv"

# Move the active selection to A2
$ws.Range("A2").Select()
